$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 5; Klasse = "05D" },
    @{ Row = 6; Klasse = "05E" },
    @{ Row = 7; Klasse = "05F" },
    @{ Row = 8; Klasse = "05G" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Klasse
    $ws.Range("B$row").Value = "Andrea Fischer, StDin"
    $ws.Range("C$row").Value = "14:30"
    $ws.Range("D$row").Value = "15:00"
    $ws.Range("E$row").Value = "15.07.2017"
}

[void]$ws.Range("B6").Select()
